$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value2 = 0.05133061381894777
$ws.Cells.Item(2, 4).Value2 = 0.1741290836617964
$ws.Cells.Item(2, 5).Value2 = 0.02897493150599217
$ws.Cells.Item(2, 6).Value2 = 9.640737515037699
$ws.Cells.Item(2, 7).Value2 = 0.002818603874210813
$ws.Cells.Item(2, 9).Value2 = 7.229883981341231
$ws.Cells.Item(2, 10).Value2 = 0.1305223454264448
$ws.Cells.Item(2, 12).Value2 = 0.06925044847813666
$ws.Cells.Item(2, 13).Value2 = 4.867458019795976
$ws.Cells.Item(2, 14).Value2 = 1.549158877228081

$ws.Cells.Item(3, 3).Value2 = 0.04592809670845099
$ws.Cells.Item(3, 4).Value2 = 0.1637124510925645
$ws.Cells.Item(3, 5).Value2 = 0.02632812220419822
$ws.Cells.Item(3, 6).Value2 = 9.581464372333869
$ws.Cells.Item(3, 7).Value2 = 0.002832094617705871
$ws.Cells.Item(3, 9).Value2 = 7.182335385552818
$ws.Cells.Item(3, 10).Value2 = 0.1323454732424327
$ws.Cells.Item(3, 12).Value2 = 0.06736256518260575
$ws.Cells.Item(3, 13).Value2 = 4.668797475599206
$ws.Cells.Item(3, 14).Value2 = 1.482098851307057

$ws.Cells.Item(4, 3).Value2 = 0.04264993749521295
$ws.Cells.Item(4, 4).Value2 = 0.1574809065674572
$ws.Cells.Item(4, 5).Value2 = 0.02469600121152382
$ws.Cells.Item(4, 6).Value2 = 9.552248600161306
$ws.Cells.Item(4, 7).Value2 = 0.00284078433586495
$ws.Cells.Item(4, 9).Value2 = 7.158563980134943
$ws.Cells.Item(4, 10).Value2 = 0.1335180159550236
$ws.Cells.Item(4, 12).Value2 = 0.06619946603236571
$ws.Cells.Item(4, 13).Value2 = 4.548852015357483
$ws.Cells.Item(4, 14).Value2 = 1.441522863053819

$ws.Cells.Item(5, 3).Value2 = 0.04132335287685862
$ws.Cells.Item(5, 4).Value2 = 0.1549815849825364
$ws.Cells.Item(5, 5).Value2 = 0.02402902709679466
$ws.Cells.Item(5, 6).Value2 = 9.54212777377677
$ws.Cells.Item(5, 7).Value2 = 0.002844428196422941
$ws.Cells.Item(5, 9).Value2 = 7.150225236400786
$ws.Cells.Item(5, 10).Value2 = 0.1340092048156922
$ws.Cells.Item(5, 12).Value2 = 0.0657244713021079
$ws.Cells.Item(5, 13).Value2 = 4.500478524012152
$ws.Cells.Item(5, 14).Value2 = 1.425140374841874

$ws.Cells.Item(6, 3).Value2 = 0.04110362083926589
$ws.Cells.Item(6, 4).Value2 = 0.1545689599881399
$ws.Cells.Item(6, 5).Value2 = 0.02391816009827963
$ws.Cells.Item(6, 6).Value2 = 9.540554438929632
$ws.Cells.Item(6, 7).Value2 = 0.002845039477270307
$ws.Cells.Item(6, 9).Value2 = 7.148921582645215
$ws.Cells.Item(6, 10).Value2 = 0.1340915739667219
$ws.Cells.Item(6, 12).Value2 = 0.0656455360811421
$ws.Cells.Item(6, 13).Value2 = 4.49247649140824
$ws.Cells.Item(6, 14).Value2 = 1.422429346350555

$ws.Cells.Item(7, 3).Value2 = 0.04263200976272685
$ws.Cells.Item(7, 4).Value2 = 0.1574470391550307
$ws.Cells.Item(7, 5).Value2 = 0.02468701389846473
$ws.Cells.Item(7, 6).Value2 = 9.552104906623072
$ws.Cells.Item(7, 7).Value2 = 0.002840833061521983
$ws.Cells.Item(7, 9).Value2 = 7.158446081878424
$ws.Cells.Item(7, 10).Value2 = 0.1335245861627703
$ws.Cells.Item(7, 12).Value2 = 0.06619306427058547
$ws.Cells.Item(7, 13).Value2 = 4.548197595382106
$ws.Cells.Item(7, 14).Value2 = 1.441301302821387

$ws.Cells.Item(8, 3).Value2 = 0.04945941328843162
$ws.Cells.Item(8, 4).Value2 = 0.1705025712046222
$ws.Cells.Item(8, 5).Value2 = 0.02806368482518096
$ws.Cells.Item(8, 6).Value2 = 9.618797950546991
$ws.Cells.Item(8, 7).Value2 = 0.002823171479011748
$ws.Cells.Item(8, 9).Value2 = 7.212354049450028
$ws.Cells.Item(8, 10).Value2 = 0.131139936139971
$ws.Cells.Item(8, 12).Value2 = 0.06860030085626789
$ws.Cells.Item(8, 13).Value2 = 4.798534199553302
$ws.Cells.Item(8, 14).Value2 = 1.525913321850993

$ws.Cells.Item(9, 3).Value2 = 0.06318192916008059
$ws.Cells.Item(9, 4).Value2 = 0.1974664869014759
$ws.Cells.Item(9, 5).Value2 = 0.03463582216302186
$ws.Cells.Item(9, 6).Value2 = 9.807492770620911
$ws.Cells.Item(9, 7).Value2 = 0.002791736010939649
$ws.Cells.Item(9, 9).Value2 = 7.361835464153046
$ws.Cells.Item(9, 10).Value2 = 0.1268848546951782
$ws.Cells.Item(9, 12).Value2 = 0.07329151199263606
$ws.Cells.Item(9, 13).Value2 = 5.305884841453945
$ws.Cells.Item(9, 14).Value2 = 1.696520859919787

$ws.Cells.Item(10, 3).Value2 = 0.07350314427165472
$ws.Cells.Item(10, 4).Value2 = 0.2181920653218299
$ws.Cells.Item(10, 5).Value2 = 0.03944285088661559
$ws.Cells.Item(10, 6).Value2 = 9.982819488141388
$ws.Cells.Item(10, 7).Value2 = 0.002770554647560047
$ws.Cells.Item(10, 9).Value2 = 7.499409980499848
$ws.Cells.Item(10, 10).Value2 = 0.1240147779986733
$ws.Cells.Item(10, 12).Value2 = 0.07672319022252339
$ws.Cells.Item(10, 13).Value2 = 5.689153375904965
$ws.Cells.Item(10, 14).Value2 = 1.824646402142747

$ws.Cells.Item(11, 3).Value2 = 0.07825857435520334
$ws.Cells.Item(11, 4).Value2 = 0.227838633859335
$ws.Cells.Item(11, 5).Value2 = 0.04162699796450653
$ws.Cells.Item(11, 6).Value2 = 10.07086409644415
$ws.Cells.Item(11, 7).Value2 = 0.002761326533624296
$ws.Cells.Item(11, 9).Value2 = 7.568265309217878
$ws.Cells.Item(11, 10).Value2 = 0.1227645840732077
$ws.Cells.Item(11, 12).Value2 = 0.07828179815438574
$ws.Cells.Item(11, 13).Value2 = 5.865913228908312
$ws.Cells.Item(11, 14).Value2 = 1.883524234034809

$ws.Cells.Item(12, 3).Value2 = 0.08006866639124155
$ws.Cells.Item(12, 4).Value2 = 0.2315245443114122
$ws.Cells.Item(12, 5).Value2 = 0.04245386412893382
$ws.Cells.Item(12, 6).Value2 = 10.10542265114879
$ws.Cells.Item(12, 7).Value2 = 0.002757890044657566
$ws.Cells.Item(12, 9).Value2 = 7.595261388259672
$ws.Cells.Item(12, 10).Value2 = 0.1222991305798953
$ws.Cells.Item(12, 12).Value2 = 0.07887169867889554
$ws.Cells.Item(12, 13).Value2 = 5.933203307169549
$ws.Cells.Item(12, 14).Value2 = 1.905903660720696

$ws.Cells.Item(13, 3).Value2 = 0.07967840587267006
$ws.Cells.Item(13, 4).Value2 = 0.2307292260406655
$ws.Cells.Item(13, 5).Value2 = 0.04227579116799163
$ws.Cells.Item(13, 6).Value2 = 10.09792525035135
$ws.Cells.Item(13, 7).Value2 = 0.002758627583827673
$ws.Cells.Item(13, 9).Value2 = 7.589405961749492
$ws.Cells.Item(13, 10).Value2 = 0.1223990199637415
$ws.Cells.Item(13, 12).Value2 = 0.07874466609243314
$ws.Cells.Item(13, 13).Value2 = 5.918695256821024
$ws.Cells.Item(13, 14).Value2 = 1.901080153946879

$ws.Cells.Item(14, 3).Value2 = 0.07840730125397499
$ws.Cells.Item(14, 4).Value2 = 0.2281412061489902
$ws.Cells.Item(14, 5).Value2 = 0.04169502832803573
$ws.Cells.Item(14, 6).Value2 = 10.07368268590903
$ws.Cells.Item(14, 7).Value2 = 0.002761042652442734
$ws.Cells.Item(14, 9).Value2 = 7.570467695334827
$ws.Cells.Item(14, 10).Value2 = 0.1227261312532342
$ws.Cells.Item(14, 12).Value2 = 0.07833033557440672
$ws.Cells.Item(14, 13).Value2 = 5.87144205963034
$ws.Cells.Item(14, 14).Value2 = 1.885363734880514

$ws.Cells.Item(15, 3).Value2 = 0.07762994555923797
$ws.Cells.Item(15, 4).Value2 = 0.2265603093682103
$ws.Cells.Item(15, 5).Value2 = 0.0413392698428261
$ws.Cells.Item(15, 6).Value2 = 10.05899282610574
$ws.Cells.Item(15, 7).Value2 = 0.002762529487289471
$ws.Cells.Item(15, 9).Value2 = 7.558988150769181
$ws.Cells.Item(15, 10).Value2 = 0.1229275339756821
$ws.Cells.Item(15, 12).Value2 = 0.07807650714840975
$ws.Cells.Item(15, 13).Value2 = 5.842544599640888
$ws.Cells.Item(15, 14).Value2 = 1.875747824926606

$ws.Cells.Item(16, 3).Value2 = 0.07319363322763195
$ws.Cells.Item(16, 4).Value2 = 0.2175661661376864
$ws.Cells.Item(16, 5).Value2 = 0.03930007223601351
$ws.Cells.Item(16, 6).Value2 = 9.977234553374046
$ws.Cells.Item(16, 7).Value2 = 0.002771165874350555
$ws.Cells.Item(16, 9).Value2 = 7.495038037325344
$ws.Cells.Item(16, 10).Value2 = 0.124097596153975
$ws.Cells.Item(16, 12).Value2 = 0.07662128551316982
$ws.Cells.Item(16, 13).Value2 = 5.677650899176655
$ws.Cells.Item(16, 14).Value2 = 1.820810390253229

$ws.Cells.Item(17, 3).Value2 = 0.07048798738077267
$ws.Cells.Item(17, 4).Value2 = 0.2121055798201894
$ws.Cells.Item(17, 5).Value2 = 0.03804853108941231
$ws.Cells.Item(17, 6).Value2 = 9.929219509931045
$ws.Cells.Item(17, 7).Value2 = 0.002776567947836183
$ws.Cells.Item(17, 9).Value2 = 7.457427228514234
$ws.Cells.Item(17, 10).Value2 = 0.1248295843512475
$ws.Cells.Item(17, 12).Value2 = 0.07572795106021601
$ws.Cells.Item(17, 13).Value2 = 5.577116897602053
$ws.Cells.Item(17, 14).Value2 = 1.787258750870308

$ws.Cells.Item(18, 3).Value2 = 0.06893740302750473
$ws.Cells.Item(18, 4).Value2 = 0.2089852552614389
$ws.Cells.Item(18, 5).Value2 = 0.03732842420339466
$ws.Cells.Item(18, 6).Value2 = 9.902380969048181
$ws.Cells.Item(18, 7).Value2 = 0.002779713459063606
$ws.Cells.Item(18, 9).Value2 = 7.436383585580046
$ws.Cells.Item(18, 10).Value2 = 0.1252558208961663
$ws.Cells.Item(18, 12).Value2 = 0.07521389239776965
$ws.Cells.Item(18, 13).Value2 = 5.519518612617787
$ws.Cells.Item(18, 14).Value2 = 1.768016629805828

$ws.Cells.Item(19, 3).Value2 = 0.06841334962884105
$ws.Cells.Item(19, 4).Value2 = 0.207932234439113
$ws.Cells.Item(19, 5).Value2 = 0.03708456036110519
$ws.Cells.Item(19, 6).Value2 = 9.893426815130624
$ws.Cells.Item(19, 7).Value2 = 0.002780785084577224
$ws.Cells.Item(19, 9).Value2 = 7.429359146015713
$ws.Cells.Item(19, 10).Value2 = 0.1254010332146605
$ws.Cells.Item(19, 12).Value2 = 0.07503979901702706
$ws.Cells.Item(19, 13).Value2 = 5.500055423213183
$ws.Cells.Item(19, 14).Value2 = 1.761511227787935

$ws.Cells.Item(20, 3).Value2 = 0.07077542052364549
$ws.Cells.Item(20, 4).Value2 = 0.2126847398775453
$ws.Cells.Item(20, 5).Value2 = 0.03818178495007629
$ws.Cells.Item(20, 6).Value2 = 9.934250040608902
$ws.Cells.Item(20, 7).Value2 = 0.002775988919817879
$ws.Cells.Item(20, 9).Value2 = 7.461369850136407
$ws.Cells.Item(20, 10).Value2 = 0.1247511231016007
$ws.Cells.Item(20, 12).Value2 = 0.075823072169932
$ws.Cells.Item(20, 13).Value2 = 5.587795441149183
$ws.Cells.Item(20, 14).Value2 = 1.790824607030345

$ws.Cells.Item(21, 3).Value2 = 0.07878039743329168
$ws.Cells.Item(21, 4).Value2 = 0.228900462990282
$ws.Cells.Item(21, 5).Value2 = 0.04186561727741989
$ws.Cells.Item(21, 6).Value2 = 10.08077004540297
$ws.Cells.Item(21, 7).Value2 = 0.00276033171893228
$ws.Cells.Item(21, 9).Value2 = 7.576005134399452
$ws.Cells.Item(21, 10).Value2 = 0.122629834501998
$ws.Cells.Item(21, 12).Value2 = 0.07845204253699478
$ws.Cells.Item(21, 13).Value2 = 5.885311765146099
$ws.Cells.Item(21, 14).Value2 = 1.88997776960727

$ws.Cells.Item(22, 3).Value2 = 0.0840667184785957
$ws.Cells.Item(22, 4).Value2 = 0.2396912797867969
$ws.Cells.Item(22, 5).Value2 = 0.04427198911603725
$ws.Cells.Item(22, 6).Value2 = 10.18364149270769
$ws.Cells.Item(22, 7).Value2 = 0.002750436610082446
$ws.Cells.Item(22, 9).Value2 = 7.656310290768289
$ws.Cells.Item(22, 10).Value2 = 0.1212898919541825
$ws.Cells.Item(22, 12).Value2 = 0.08016844024440672
$ws.Cells.Item(22, 13).Value2 = 6.081829164461482
$ws.Cells.Item(22, 14).Value2 = 1.955267103695775

$ws.Cells.Item(23, 3).Value2 = 0.0812400926568273
$ws.Cells.Item(23, 4).Value2 = 0.2339138375491245
$ws.Cells.Item(23, 5).Value2 = 0.04298772418426466
$ws.Cells.Item(23, 6).Value2 = 10.12807738033081
$ws.Cells.Item(23, 7).Value2 = 0.002755687104708903
$ws.Cells.Item(23, 9).Value2 = 7.612950390047246
$ws.Cells.Item(23, 10).Value2 = 0.122000795376632
$ws.Cells.Item(23, 12).Value2 = 0.07925251402679123
$ws.Cells.Item(23, 13).Value2 = 5.976751475580272
$ws.Cells.Item(23, 14).Value2 = 1.920376899301942

$ws.Cells.Item(24, 3).Value2 = 0.07064545669938127
$ws.Cells.Item(24, 4).Value2 = 0.2124228424299588
$ws.Cells.Item(24, 5).Value2 = 0.0381215426951016
$ws.Cells.Item(24, 6).Value2 = 9.931973353907438
$ws.Cells.Item(24, 7).Value2 = 0.002776250574247961
$ws.Cells.Item(24, 9).Value2 = 7.459585587989409
$ws.Cells.Item(24, 10).Value2 = 0.1247865785915954
$ws.Cells.Item(24, 12).Value2 = 0.07578006937988846
$ws.Cells.Item(24, 13).Value2 = 5.582967049134169
$ws.Cells.Item(24, 14).Value2 = 1.789212336435497

$ws.Cells.Item(25, 3).Value2 = 0.05943024994589052
$ws.Cells.Item(25, 4).Value2 = 0.190017612705617
$ws.Cells.Item(25, 5).Value2 = 0.03286230717874261
$ws.Cells.Item(25, 6).Value2 = 9.750105567728156
$ws.Cells.Item(25, 7).Value2 = 0.002799901359779626
$ws.Cells.Item(25, 9).Value2 = 7.316602978856196
$ws.Cells.Item(25, 10).Value2 = 0.1279909418732919
$ws.Cells.Item(25, 12).Value2 = 0.07202526236003592
$ws.Cells.Item(25, 13).Value2 = 5.166823247771759
$ws.Cells.Item(25, 14).Value2 = 1.649873639787614
